# Update "Luy ke thang SOC TRANG" rows with refreshed Notion export values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2024-07-24T16:01:00.000Z"

# Row 2
$ws.Range("D2").Value = $newTimestamp
$ws.Range("W2").Value = 96219000
$ws.Range("AA2").Value = 23498000
$ws.Range("AE2").Value = 119895000
$ws.Range("AH2").Value = 106895000

# Row 3
$ws.Range("D3").Value = $newTimestamp
$ws.Range("W3").Value = 141611000
$ws.Range("AA3").Value = 74400000
$ws.Range("AE3").Value = 470450000
$ws.Range("AH3").Value = 447050000

# Row 6
$ws.Range("D6").Value = $newTimestamp
$ws.Range("W6").Value = 227830000
$ws.Range("AA6").Value = 28800000
$ws.Range("AE6").Value = 402100000
$ws.Range("AH6").Value = 382100000

# Row 8
$ws.Range("D8").Value = $newTimestamp
$ws.Range("W8").Value = 112842000
$ws.Range("AA8").Value = 44286000
$ws.Range("AE8").Value = 133986000
$ws.Range("AH8").Value = 114800000

# Row 11
$ws.Range("D11").Value = $newTimestamp
$ws.Range("W11").Value = 200504000
$ws.Range("AA11").Value = 59600000
$ws.Range("AE11").Value = 306700000
$ws.Range("AH11").Value = 244700000

# Row 13
$ws.Range("D13").Value = $newTimestamp
$ws.Range("W13").Value = 54043000
$ws.Range("AA13").Value = 15000000
$ws.Range("AE13").Value = 205588000
$ws.Range("AH13").Value = 179588000
$ws.Range("AK13").Value = 20
$ws.Range("AQ13").Value = 194588000
